# Updates the handback-status report with a new CI run:
#   - old handoff UUID 4119853f-5993-4abe-b6dd-b7b10fafb661 -> a76d7af6-a1e2-4e5f-92c8-d5eafea6df61
#   - old handoff UUID f44aacc3-2cbe-4fb0-8022-135d9610169f -> ffff0bfc01cc-7c5e-41ce-a979-fa12fabdc253
#   - old content-hash ca095af50ed17dd4391394173bfdb6138a2a036c -> 6429452d5d1da290753ce6a1c9916897a8eefebd
#   - new timestamps for the latest HO Xliff generation / handback round-trip
#
# The workbook has three sheets:
#   "Overview" - summary table (File Name / Path And Name / ... / Latest HO Xliff Generate Date)
#   "zh-cn"    - per-locale detail table for zh-cn
#   "de-de"    - per-locale detail table for de-de

$wb = $excel.ActiveWorkbook

$oldId1 = "4119853f-5993-4abe-b6dd-b7b10fafb661"
$newId1 = "a76d7af6-a1e2-4e5f-92c8-d5eafea6df61"
$oldId2 = "f44aacc3-2cbe-4fb0-8022-135d9610169f"
$newId2 = "ffff0bfc01cc-7c5e-41ce-a979-fa12fabdc253"

$newXlfZhCn = "$newId1.6429452d5d1da290753ce6a1c9916897a8eefebd.zh-cn.xlf"
$newXlfDeDe = "$newId1.6429452d5d1da290753ce6a1c9916897a8eefebd.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("B2").Value = "e2e\$newId1.md"
$ws.Range("G2").Value = "2016-08-17 00:58:40"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("B3").Value = "e2e\$newId2.md"
$ws.Range("G3").Value = "2016-08-17 00:58:40"

$hyperlinks = @($ws.Hyperlinks)
foreach ($hl in $hyperlinks) {
    $addr = $hl.Address
    if ($addr -like "*$oldId1*") {
        $hl.TextToDisplay = "e2e\$newId1.md"
    } elseif ($addr -like "*$oldId2*") {
        $hl.TextToDisplay = "e2e\$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("G2").Value = $newXlfZhCn
$ws.Range("H2").Value = "2016-08-17 00:58:35"
$ws.Range("I2").Value = "$newId1.md"
$ws.Range("J2").Value = $newXlfZhCn
$ws.Range("K2").Value = "2016-08-17 00:58:51"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("G3").Value = $newXlfZhCn
$ws.Range("H3").Value = "2016-08-17 00:58:35"
$ws.Range("I3").Value = "$newId2.md"
$ws.Range("J3").Value = $newXlfZhCn
$ws.Range("K3").Value = "2016-08-17 00:58:51"

$hyperlinks = @($ws.Hyperlinks)
foreach ($hl in $hyperlinks) {
    $addr = $hl.Address
    if ($addr -like "*$oldId1*") {
        $hl.TextToDisplay = "$newId1.md"
    } elseif ($addr -like "*$oldId2*") {
        $hl.TextToDisplay = "$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("G2").Value = $newXlfDeDe
$ws.Range("H2").Value = "2016-08-17 00:58:40"
$ws.Range("I2").Value = "$newId1.md"
$ws.Range("J2").Value = $newXlfDeDe
$ws.Range("K2").Value = "2016-08-17 00:58:58"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("G3").Value = $newXlfDeDe
$ws.Range("H3").Value = "2016-08-17 00:58:40"
$ws.Range("I3").Value = "$newId2.md"
$ws.Range("J3").Value = $newXlfDeDe
$ws.Range("K3").Value = "2016-08-17 00:58:58"

$hyperlinks = @($ws.Hyperlinks)
foreach ($hl in $hyperlinks) {
    $addr = $hl.Address
    if ($addr -like "*$oldId1*") {
        $hl.TextToDisplay = "$newId1.md"
    } elseif ($addr -like "*$oldId2*") {
        $hl.TextToDisplay = "$newId2.md"
    }
}
